$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 51064.8
$ws.Range("I70").Value = 2083.3333
$ws.Range("K70").Value = 6249.999899999999
$ws.Range("M70").Value = -5979.999899999999

$ws.Range("H73").Value = 51064.8
$ws.Range("I73").Value = 2083.3333
$ws.Range("K73").Value = 6249.999899999999
$ws.Range("M73").Value = -5313.999899999999

$ws.Range("H112").Value = 1261.3334
$ws.Range("J112").Value = 1454.5
$ws.Range("L112").Value = 4363.5
$ws.Range("N112").Value = -6579.5

$ws.Range("H137").Value = 1517
$ws.Range("I137").Value = 1310
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 3930
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -1380
$ws.Range("N137").Value = -11100

$ws.Range("H138").Value = 3893.3936
$ws.Range("I138").Value = 3274.4546
$ws.Range("J138").Value = 4622.857
$ws.Range("K138").Value = 9823.363799999999
$ws.Range("L138").Value = 13868.571
$ws.Range("M138").Value = -4683.363799999999
$ws.Range("N138").Value = -24148.571

$ws.Range("H141").Value = 1901.7931
$ws.Range("I141").Value = 1434
$ws.Range("J141").Value = 15000
$ws.Range("K141").Value = 4302
$ws.Range("L141").Value = 45000
$ws.Range("M141").Value = 878
$ws.Range("N141").Value = -55360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5790.9736
$ws.Range("I32").Value = 3335.3635
$ws.Range("J32").Value = 21998
$ws.Range("K32").Value = 3335.3635
$ws.Range("L32").Value = 21998
$ws.Range("M32").Value = -3048.3635
$ws.Range("N32").Value = -22572

$ws.Range("H36").Value = 11506.5
$ws.Range("I36").Value = 8675.333000000001
$ws.Range("K36").Value = 8675.333000000001
$ws.Range("M36").Value = -8329.333000000001

$ws.Range("H45").Value = 4499.5713
$ws.Range("I45").Value = 4499.5713
$ws.Range("K45").Value = 4499.5713
$ws.Range("M45").Value = -4122.5713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2035.6
$ws.Range("I20").Value = 1939.125
$ws.Range("K20").Value = 1939.125
$ws.Range("M20").Value = -1692.125

$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372

$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864

$ws.Range("H86").Value = 1667.909
$ws.Range("I86").Value = 1457.6666
$ws.Range("J86").Value = 1920.2
$ws.Range("K86").Value = 1457.6666
$ws.Range("L86").Value = 1920.2
$ws.Range("M86").Value = -334.6666
$ws.Range("N86").Value = -4166.2

$ws.Range("H89").Value = 1667.909
$ws.Range("I89").Value = 1457.6666
$ws.Range("J89").Value = 1920.2
$ws.Range("K89").Value = 7288.333000000001
$ws.Range("L89").Value = 9601
$ws.Range("M89").Value = -1672.333000000001
$ws.Range("N89").Value = -20833

$ws.Range("H100").Value = 42643
$ws.Range("J100").Value = 42643
$ws.Range("L100").Value = 42643
$ws.Range("N100").Value = -44807

$ws.Range("H134").Value = 3143.5217
$ws.Range("I134").Value = 2342.111
$ws.Range("K134").Value = 7026.333
$ws.Range("M134").Value = -4491.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 418.07144
$ws.Range("I22").Value = 261.75
$ws.Range("K22").Value = 261.75
$ws.Range("M22").Value = 88.25

$ws.Range("H31").Value = 3812.5454
$ws.Range("I31").Value = 3348.111
$ws.Range("K31").Value = 3348.111
$ws.Range("M31").Value = -3053.111

$ws.Range("H34").Value = 3812.5454
$ws.Range("I34").Value = 3348.111
$ws.Range("K34").Value = 3348.111
$ws.Range("M34").Value = -3146.111

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").Value = ""

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = ""

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = ""

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""

$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37372

$ws.Range("H88").Value = 92822.75
$ws.Range("J88").Value = 92822.75
$ws.Range("L88").Value = 92822.75
$ws.Range("N88").Value = -93634.75

$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -116856

$ws.Range("H91").Value = 92822.75
$ws.Range("J91").Value = 92822.75
$ws.Range("L91").Value = 92822.75
$ws.Range("N91").Value = -95630.75

$ws.Range("H99").Value = 12722.36
$ws.Range("I99").Value = 9252.454
$ws.Range("J99").Value = 15448.714
$ws.Range("K99").Value = 9252.454
$ws.Range("L99").Value = 15448.714
$ws.Range("M99").Value = -7754.454
$ws.Range("N99").Value = -18444.714

$ws.Range("H107").Value = 1128.5454
$ws.Range("I107").Value = 732.2857
$ws.Range("J107").Value = 1822
$ws.Range("K107").Value = 732.2857
$ws.Range("L107").Value = 1822
$ws.Range("M107").Value = 1187.7143
$ws.Range("N107").Value = -5662

$ws.Range("H126").Value = 12722.36
$ws.Range("I126").Value = 9252.454
$ws.Range("J126").Value = 15448.714
$ws.Range("K126").Value = 27757.362
$ws.Range("L126").Value = 46346.142
$ws.Range("M126").Value = -25287.362
$ws.Range("N126").Value = -51286.142

$ws.Range("H132").Value = 1421.7826
$ws.Range("I132").Value = 1440.1
$ws.Range("J132").Value = 1299.6666
$ws.Range("K132").Value = 4320.299999999999
$ws.Range("L132").Value = 3898.9998
$ws.Range("M132").Value = -1790.299999999999
$ws.Range("N132").Value = -8958.9998

$ws.Range("H134").Value = 2091.6216
$ws.Range("I134").Value = 1891.931
$ws.Range("K134").Value = 5675.793
$ws.Range("M134").Value = -3140.793

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 478.57144
$ws.Range("J2").Value = 1100.3334
$ws.Range("L2").Value = 1100.3334
$ws.Range("N2").Value = -1326.3334

$ws.Range("H93").Value = 59326.668
$ws.Range("J93").Value = 59326.668
$ws.Range("L93").Value = 59326.668
$ws.Range("N93").Value = -63070.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4890.5
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = ""

$ws.Range("H16").Value = 3001.5
$ws.Range("J16").Value = 3001.5
$ws.Range("L16").Value = 3001.5
$ws.Range("N16").Value = -3341.5

$ws.Range("H22").Value = 3947.2273
$ws.Range("I22").Value = 2846.7896
$ws.Range("J22").Value = 10916.667
$ws.Range("K22").Value = 2846.7896
$ws.Range("L22").Value = 10916.667
$ws.Range("M22").Value = -2551.7896
$ws.Range("N22").Value = -11506.667

$ws.Range("H27").Value = 3947.2273
$ws.Range("I27").Value = 2846.7896
$ws.Range("J27").Value = 10916.667
$ws.Range("K27").Value = 2846.7896
$ws.Range("L27").Value = 10916.667
$ws.Range("M27").Value = -2739.7896
$ws.Range("N27").Value = -11130.667

$ws.Range("H40").Value = 2837.5386
$ws.Range("I40").Value = 2963.2727
$ws.Range("J40").Value = 2146
$ws.Range("K40").Value = 2963.2727
$ws.Range("L40").Value = 2146
$ws.Range("M40").Value = -2827.2727
$ws.Range("N40").Value = -2418

$ws.Range("H46").Value = 4765
$ws.Range("I46").Value = 1799
$ws.Range("K46").Value = 1799
$ws.Range("M46").Value = -1611

$ws.Range("H95").Value = 27200
$ws.Range("J95").Value = 27200
$ws.Range("L95").Value = 27200
$ws.Range("N95").Value = -32692

$ws.Range("H126").Value = 4890.5
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = ""

$ws.Range("H132").Value = 4341.75
$ws.Range("I132").Value = 3775.5
$ws.Range("K132").Value = 11326.5
$ws.Range("M132").Value = -8796.5

$ws.Range("H136").Value = 2984.3914
$ws.Range("I136").Value = 2464.4285
$ws.Range("K136").Value = 7393.2855
$ws.Range("M136").Value = -4843.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 23699.75
$ws.Range("I41").Value = 50000
$ws.Range("J41").Value = 14933
$ws.Range("K41").Value = 50000
$ws.Range("L41").Value = 14933
$ws.Range("M41").Value = -49610
$ws.Range("N41").Value = -15713

$ws.Range("H45").Value = 15789
$ws.Range("J45").Value = 16188
$ws.Range("L45").Value = 16188
$ws.Range("N45").Value = -17170

$ws.Range("H97").Value = 41572
$ws.Range("J97").Value = 41572
$ws.Range("L97").Value = 41572
$ws.Range("N97").Value = -43554
